# Update gh-pages to output generated at 456a3b4
# This script updates the "想去人数" (number of people interested) figures
# on the "展览" (Exhibition) and "全部类型" (All Categories) sheets to
# reflect freshly-scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 19
$wsExhibit.Range("F7").Value = 419
$wsExhibit.Range("F10").Value = 513

# --- Sheet "全部类型" (All Categories) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 0
$wsAll.Range("F8").Value = 147
$wsAll.Range("F10").Value = 513
